# Append a new daily-log row (row 35) to each of the 4 data sheets,
# mirroring the previous row (row 34) except for the updated timestamp
# (column A) and, for sheet 1 only, the updated D/H values.
#
# This corresponds to the upstream diff that grows each sheet's
# <dimension> from A1:I34 to A1:I35 and appends one fully-populated row.

function Add-LogRow {
    param($ws, $row, $timeValue, $bVal, $cVal, $dVal, $eVal, $fVal, $gVal, $hVal, $iVal)

    $prevRow = $row - 1

    # Column A: date/time serial, formatted the same as the row above it.
    $ws.Range("A$row").Value2 = $timeValue
    $ws.Range("A$row").NumberFormat = $ws.Range("A$prevRow").NumberFormat

    # Columns B-E: hex-string payloads (kept as plain text).
    $ws.Range("B$row").Value = $bVal
    $ws.Range("C$row").Value = $cVal
    $ws.Range("D$row").Value = $dVal
    $ws.Range("E$row").Value = $eVal

    # Columns F-I: decimal counterparts.
    $ws.Range("F$row").Value = $fVal
    $ws.Range("G$row").Value = $gVal
    $ws.Range("H$row").Value = $hVal
    $ws.Range("I$row").Value = $iVal
}

$wb = $excel.ActiveWorkbook

# The new row's timestamp (same instant across all 4 sheets).
$newTime = 45821.4962037037

# --- Sheet 1: FE_LFT_#1 ---
$ws1 = $wb.Worksheets.Item(1)
$g1 = [double]"7.598631275147109e+23"
Add-LogRow $ws1 35 $newTime "0x01,0x7c" "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0," "0x01,0x6C" "0xf" 380 $g1 364 15

# --- Sheet 2: FE_LFT_#2 ---
$ws2 = $wb.Worksheets.Item(2)
$g2 = [double]"5.68432987514711e+23"
Add-LogRow $ws2 35 $newTime "0x01,0x90" "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78," "0x01,0x80" "0xe" 400 $g2 384 14

# --- Sheet 3: FE_PLT_#1 ---
$ws3 = $wb.Worksheets.Item(3)
$g3 = [double]"5.68631262647114e+23"
Add-LogRow $ws3 35 $newTime "0x00,0x6e" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c," "0x00,0x6C" "0x3" 110 $g3 108 3

# --- Sheet 4: FE_PLT_#2 ---
$ws4 = $wb.Worksheets.Item(4)
$g4 = [double]"9.85046333984776e+23"
Add-LogRow $ws4 35 $newTime "0x00,0x6e" "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c," "0x00,0x6C" "0x3" 110 $g4 108 3

Write-Host "Appended row 35 to all 4 sheets"
